# Updates the workbook for data through 2022-07-20 (adds 2022-07-28 data per commit msg,
# but per the diff it actually updates the "through July 19" -> "through July 20" counts).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label for the current month column.
$ws.Name = "Through 2022-07-20"
$ws.Range("B1").Value = "July 2022 (through July 20)"

# Apply the updated/added cell values (row => column => new value)
$updates = @{
    2  = @{ "B" = 12; "AD" = 9 }
    3  = @{ "AY" = 1 }
    4  = @{ "AK" = 2 }
    5  = @{ "B" = 7 }
    8  = @{ "P" = 13 }
    12 = @{ "AD" = 1 }
    15 = @{ "I" = 1; "AY" = 1 }
    16 = @{ "B" = 2 }
    20 = @{ "I" = 6; "P" = 4; "AD" = 2 }
    27 = @{ "B" = 5 }
    29 = @{ "B" = 5; "AY" = 2 }
    38 = @{ "AR" = 2 }
    39 = @{ "AD" = 1 }
    41 = @{ "B" = 2; "W" = 1 }
    49 = @{ "W" = 1 }
    50 = @{ "AR" = 1 }
    52 = @{ "I" = 5 }
    56 = @{ "AK" = 2 }
    61 = @{ "AR" = 1 }
    65 = @{ "AK" = 1 }
    79 = @{ "I" = 2 }
    94 = @{ "P" = 2 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
